$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("B1").Value = "Furo"
$ws.Range("C1").Value = "Pino"
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null

# Update existing data rows
$ws.Range("B2").Value = 40.20018007729973
$ws.Range("C2").Value = 34.12010804637984

$ws.Range("B3").Value = 59.79981992270027
$ws.Range("C3").Value = 34.12010804637984

# Add new rows 4 and 5
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 40.20018007729973
$ws.Range("C4").Value = 45.87989195362016

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 59.79981992270027
$ws.Range("C5").Value = 45.87989195362016

# Copy the style from A2/A3 (bold, centered, bordered) to A4 and A5
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null
